$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;11;14.json"
$ws.Range("C2").Value = "MSG: None`n`nMSG: The theme for spirit week has been decided through a compromise.`n"
$ws.Range("D2").Value = "they_came_to_a_compromise, "

$ws.Range("B3").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;24;22.json"
$ws.Range("C3").Value = "MSG: None`n`nMSG: There was no consensus reached on the theme for spirit week.`n"

$ws.Range("B4").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;29;20.json"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The selected theme for spirit week is `"Crazy Sock and Hat Day.`"`n"
$ws.Range("D4").Value = "they_came_to_a_compromise, "

$ws.Range("B5").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;07;52.json"
$ws.Range("C5").Value = "MSG: None`n`nMSG: Based on the discussion, it seems there was no consensus reached on selecting a theme for spirit week.`n"

$ws.Range("B6").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;09;43.json"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision for the spirit week theme resulted in no consensus.`n"

$ws.Range("B7").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;28;39.json"

$ws.Range("B8").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;06;17.json"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The outcome of the discussion resulted in no consensus for the theme of spirit week.`n"

$ws.Range("B9").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;17;01.json"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The theme for spirit week has been selected based on the consensus reached by both parties. They agreed on a combined theme, which is likely a fun and collaborative idea.`n"
$ws.Range("D9").Value = "they_came_to_a_compromise, "

$ws.Range("B10").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;26;20.json"
$ws.Range("C10").Value = "MSG: None`n`nMSG: There was no consensus reached on the theme for spirit week.`n"

$ws.Range("B11").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;09;14.json"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The discussion did not result in a consensus on the theme for spirit week.`n"

$ws.Range("B12").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;21;44.json"
$ws.Range("C12").Value = "MSG: None`n`nMSG: There was no consensus on the theme for spirit week, so the decision has been recorded.`n"

$ws.Range("B13").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;14;32.json"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The selected theme for spirit week is a combination of Crazy Sock Day and Crazy Hat Day, as both parties came to a compromise.`n"
$ws.Range("D13").Value = "they_came_to_a_compromise, "

$ws.Range("B14").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;25;03.json"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The theme for spirit week is a combination of Crazy Hat Day and Crazy Sock Day, as both parties came to a compromise.`n"

$ws.Range("B15").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;13;01.json"
$ws.Range("C15").Value = "MSG: None`n`nMSG: There was no consensus reached on the theme for spirit week, so the appropriate result has been called.`n"

$ws.Range("B16").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;20;50.json"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision resulted in no consensus on the theme for spirit week.`n"

$ws.Range("B17").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;13;51.json"
$ws.Range("C17").Value = "MSG: None`n`nMSG: There was no consensus reached on the theme for spirit week.`n"

$ws.Range("B18").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;12;17.json"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The outcome of the discussion resulted in no consensus.`n"
$ws.Range("D18").Value = "no_consensus, "

$ws.Range("B19").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;25;34.json"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The outcome is that there was no consensus reached on the theme for spirit week.`n"

$ws.Range("B20").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;23;30.json"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The discussion did not result in a consensus, so I have called the function for no consensus.`n"
$ws.Range("D20").Value = "no_consensus, "

$ws.Range("B21").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;08;38.json"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The theme for spirit week has been decided as Crazy Sock Day, following a compromise between both parties.`n"

$ws.Range("B22").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;07;15.json"
$ws.Range("C22").Value = "MSG: None`n`nMSG: There was no consensus reached on the theme for spirit week.`n"

$ws.Range("B23").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;27;35.json"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The outcome of the discussion resulted in no consensus.`n"

$ws.Range("B24").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;16;35.json"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision for spirit week has been made: a hybrid theme combining both Crazy Hat Day and Crazy Sock Day has been selected!`n"
$ws.Range("D24").Value = "they_came_to_a_compromise, "

$ws.Range("B25").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;10;31.json"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The outcome of the discussion is that there was no consensus reached on a theme for spirit week.`n"

$ws.Range("B26").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;19;59.json"
$ws.Range("C26").Value = "MSG: None`n`nMSG: There was no consensus reached regarding the theme for spirit week.`n"

$ws.Range("B27").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;16;03.json"
$ws.Range("C27").Value = "MSG: None`n`nMSG: It seems there was no consensus reached for the theme of spirit week.`n"

$ws.Range("B28").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;21;15.json"
$ws.Range("C28").Value = "MSG: None`n`nMSG: There was no consensus reached for the theme for spirit week.`n"

$ws.Range("B29").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_11;15;01.json"
$ws.Range("C29").Value = "MSG: None`n`nMSG: There was no consensus reached on the theme for spirit week.`n"

$ws.Range("B30").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;22;37.json"
$ws.Range("C30").Value = "MSG: None`n`nMSG: Since there was no consensus reached regarding the theme for spirit week, I have called the function for no consensus.`n"

$ws.Range("B31").Value = "./Warehouse/BA/run4o_discovery_03_28_2025 at_15;19;20.json"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The theme for spirit week is Crazy Sock Day, as decided by both parties.`n"
$ws.Range("D31").Value = "crazy_sock_day_was_selected, "
